$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

# Remove existing hyperlinks (F2:F13) before clearing the rows, otherwise the
# <hyperlinks>/relationship entries are left dangling.
$ws.Range("F2:F13").Hyperlinks.Delete()

# Drop all the old data rows (2-13); only the header row (1) stays.
$ws.Rows("2:13").Delete()

# --- Row 2 ---
$ws.Cells.Item(2,1).Value = "2025-09-30 06:28:04"
$ws.Cells.Item(2,2).Value = "【急募】LINE WORKSで定期メッセージ配信ツール作成依頼"
$ws.Cells.Item(2,3).Value = "システム開発"
$ws.Cells.Item(2,4).Value = "5,000 円 ~ 10,000 円 / 固定"
$ws.Cells.Item(2,5).Value = "期限情報なし"
$ws.Cells.Item(2,6).Value = "https://www.lancers.jp/work/detail/5403166"
$ws.Cells.Item(2,7).Value = 65
$ws.Cells.Item(2,8).Value = "◆ツール"

# --- Row 3 ---
$ws.Cells.Item(3,1).Value = "2025-09-30 06:28:04"
$ws.Cells.Item(3,2).Value = "【急募】教育系のWEBサイトの作成"
$ws.Cells.Item(3,3).Value = "システム開発"
$ws.Cells.Item(3,4).Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Cells.Item(3,5).Value = "期限情報なし"
$ws.Cells.Item(3,6).Value = "https://www.lancers.jp/work/detail/5403527"
$ws.Cells.Item(3,7).Value = 33
$ws.Cells.Item(3,8).Value = "◇サイト"

# --- Row 4 ---
$ws.Cells.Item(4,1).Value = "2025-09-30 06:28:04"
$ws.Cells.Item(4,2).Value = "【ビジネスパートナー募集】英語が得意な方を探しています"
$ws.Cells.Item(4,3).Value = "システム開発"
$ws.Cells.Item(4,4).Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Cells.Item(4,5).Value = "期限情報なし"
$ws.Cells.Item(4,6).Value = "https://www.lancers.jp/work/detail/5403384"
$ws.Cells.Item(4,7).Value = 25

# --- Row 5 ---
$ws.Cells.Item(5,1).Value = "2025-09-30 06:28:04"
$ws.Cells.Item(5,2).Value = "初回 サブスクペイからCSVデータをダウンロードし、データベース同期するプログラムの作成"
$ws.Cells.Item(5,3).Value = "システム開発"
$ws.Cells.Item(5,4).Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Cells.Item(5,5).Value = "期限情報なし"
$ws.Cells.Item(5,6).Value = "https://www.lancers.jp/work/detail/5403072"
$ws.Cells.Item(5,7).Value = 18

# --- Row 6 ---
$ws.Cells.Item(6,1).Value = "2025-09-30 06:28:04"
$ws.Cells.Item(6,2).Value = "限定公開 PR 限定公開の仕事"
$ws.Cells.Item(6,3).Value = "システム開発"
$ws.Cells.Item(6,4).Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Cells.Item(6,5).Value = "期限情報なし"
$ws.Cells.Item(6,6).Value = "https://www.lancers.jp/work/detail/5399347"
$ws.Cells.Item(6,7).Value = 13

# Hyperlinks + styling for the URL column
$urls = @(
    "https://www.lancers.jp/work/detail/5403166",
    "https://www.lancers.jp/work/detail/5403527",
    "https://www.lancers.jp/work/detail/5403384",
    "https://www.lancers.jp/work/detail/5403072",
    "https://www.lancers.jp/work/detail/5399347"
)
for ($r = 2; $r -le 6; $r++) {
    $cell = $ws.Cells.Item($r, 6)
    $ws.Hyperlinks.Add($cell, $urls[$r - 2]) | Out-Null
    $cell.Style = "Hyperlink"
}
